$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Rana Abo-Zaid, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy"
$ws.Cells.Item(4, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Hend Mahmoud"
$ws.Cells.Item(5, 7).Value = "Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(6, 7).Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(8, 7).Value = "Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(9, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(10, 7).Value = "Dr. Amany Raafat, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Marina Youhanna"
$ws.Cells.Item(11, 7).Value = "Dr. Sarah Mahdy, Dr. Alaa Ashraf"
$ws.Cells.Item(15, 7).Value = "Dr. Hana Amr, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(18, 7).Value = "Dr. Aya Hanafy, Dr. Yasmin, Dr. Remon, Dr. Shorok Mohammad"
$ws.Cells.Item(19, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(20, 7).Value = "Dr. Marina Sorial, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(21, 7).Value = "Dr. Yassmen Ahmad, Dr. Yasmin, Dr. Monica, Dr. Shorok Mohammad"
$ws.Cells.Item(22, 7).Value = "Dr. Naema Gomaa, Dr. Monica, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(23, 7).Value = "Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Cells.Item(24, 7).Value = "Dr. Rana Abo-Zaid, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy"
$ws.Cells.Item(25, 7).Value = "Dr. Gehan Adel, Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser"
$ws.Cells.Item(26, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Hend Mahmoud"
$ws.Cells.Item(27, 7).Value = "Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(28, 7).Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(30, 7).Value = "Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(32, 7).Value = "Dr. Amany Raafat, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Marina Youhanna"
$ws.Cells.Item(33, 7).Value = "Dr. Sarah Mahdy, Dr. Alaa Ashraf"
$ws.Cells.Item(37, 7).Value = "Dr. Hana Amr, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(40, 7).Value = "Dr. Aya Hanafy, Dr. Yasmin, Dr. Remon, Dr. Shorok Mohammad"
$ws.Cells.Item(41, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(42, 7).Value = "Dr. Marina Sorial, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(43, 7).Value = "Dr. Yassmen Ahmad, Dr. Yasmin, Dr. Monica, Dr. Shorok Mohammad"
$ws.Cells.Item(44, 7).Value = "Dr. Naema Gomaa, Dr. Monica, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(45, 7).Value = "Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Cells.Item(46, 7).Value = "Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub"
$ws.Cells.Item(48, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Cells.Item(49, 7).Value = "Dr. Amera Ahmad Saad, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(52, 7).Value = "Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din"
$ws.Cells.Item(54, 7).Value = "Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa"
$ws.Cells.Item(62, 7).Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
$ws.Cells.Item(63, 7).Value = "Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(65, 7).Value = "Dr. Shorok Mohammad, Dr. Salma Hassan, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy"
$ws.Cells.Item(66, 7).Value = "Dr. Maryam Ashraf, Dr. Marina Sorial, Dr. Monica, Dr. Aya Hanafy, Dr. Eman Mohammad Al"
$ws.Cells.Item(68, 7).Value = "Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub"
$ws.Cells.Item(70, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Cells.Item(71, 7).Value = "Dr. Amera Ahmad Saad, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(74, 7).Value = "Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din"
$ws.Cells.Item(75, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(76, 7).Value = "Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa"
$ws.Cells.Item(84, 7).Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
$ws.Cells.Item(85, 7).Value = "Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(87, 7).Value = "Dr. Shorok Mohammad, Dr. Salma Hassan, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy"
$ws.Cells.Item(88, 7).Value = "Dr. Maryam Ashraf, Dr. Marina Sorial, Dr. Monica, Dr. Aya Hanafy, Dr. Eman Mohammad Al"
$ws.Cells.Item(90, 7).Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser"
$ws.Cells.Item(92, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Cells.Item(93, 7).Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(96, 7).Value = "Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Amal Awwad, Dr. Mariam Nour El-Din"
$ws.Cells.Item(97, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(98, 7).Value = "Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa"
$ws.Cells.Item(106, 7).Value = "Dr. Neveen Nashaat, Dr. Remon, Dr. Youstina Magdy, Dr. Monica, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(107, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Monica, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(108, 7).Value = "Dr. Marina Sorial, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(110, 7).Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(111, 7).Value = "Dr. Marina Atef, Dr. Monica, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
$ws.Cells.Item(112, 7).Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser"
$ws.Cells.Item(114, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub"
$ws.Cells.Item(115, 7).Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(118, 7).Value = "Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Amal Awwad, Dr. Mariam Nour El-Din"
$ws.Cells.Item(119, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(120, 7).Value = "Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Merna Said, Dr. Marina Youhanna, Dr. Mai Mustafa"
$ws.Cells.Item(128, 7).Value = "Dr. Neveen Nashaat, Dr. Remon, Dr. Youstina Magdy, Dr. Monica, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(129, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Monica, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(130, 7).Value = "Dr. Marina Sorial, Dr. Remon, Dr. Youstina Magdy, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida"
$ws.Cells.Item(132, 7).Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(133, 7).Value = "Dr. Marina Atef, Dr. Monica, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
$ws.Cells.Item(134, 7).Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Majorelle Magdy"
$ws.Cells.Item(137, 7).Value = "Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(140, 7).Value = "Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(141, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(142, 7).Value = "Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Merna Said"
$ws.Cells.Item(144, 7).Value = "Dr. Khadija Osama, Dr. Mona Ibrahim Hussein"
$ws.Cells.Item(150, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(151, 7).Value = "Dr. Yassmen Ahmad, Dr. Marina Atef, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(152, 7).Value = "Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Cells.Item(153, 7).Value = "Dr. Maryam Ashraf, Dr. Marina Sorial, Dr. Monica, Dr. Aya Hanafy, Dr. Eman Mohammad Al"
$ws.Cells.Item(154, 7).Value = "Dr. Salma Hassan, Dr. Naema Gomaa, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(155, 7).Value = "Dr. Marina Atef, Dr. Monica, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
$ws.Cells.Item(156, 7).Value = "Dr. Mohammad El-Tanany, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Manar Montaser"
$ws.Cells.Item(159, 7).Value = "Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Cells.Item(162, 7).Value = "Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Cells.Item(163, 7).Value = "Dr. Nourhan Mohammad, Dr. Aya Essam"
$ws.Cells.Item(164, 7).Value = "Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Merna Said"
$ws.Cells.Item(166, 7).Value = "Dr. Khadija Osama, Dr. Mona Ibrahim Hussein"
$ws.Cells.Item(172, 7).Value = "Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Nardine, Dr. Wafaa Ebida"
$ws.Cells.Item(173, 7).Value = "Dr. Yassmen Ahmad, Dr. Marina Atef, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(174, 7).Value = "Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Cells.Item(175, 7).Value = "Dr. Maryam Ashraf, Dr. Marina Sorial, Dr. Monica, Dr. Aya Hanafy, Dr. Eman Mohammad Al"
$ws.Cells.Item(176, 7).Value = "Dr. Salma Hassan, Dr. Naema Gomaa, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(177, 7).Value = "Dr. Marina Atef, Dr. Monica, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Naema Gomaa"
